# Apply updated cryptocurrency price/volume figures (Price column D, Volume(1h) column E).
# Source cells are plain text (e.g. "42.065.60", "  -2.18%  "), so force the Price
# column to Text format before assignment -- otherwise COM auto-coerces numeric-looking
# strings (like "0.490" or "17.30") into numbers and trailing zeros are lost.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.115.05"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.94"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.38"
$ws.Range("E5").Value = "  -2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.07"
$ws.Range("E6").Value = "  -7.00%  "

$ws.Range("E7").Value = "  -2.91%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -3.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.04"
$ws.Range("E10").Value = "  -5.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.99"
$ws.Range("E12").Value = "  -8.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.112"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.65"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.616.80"
$ws.Range("E15").Value = "  -3.27%  "

$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.275.61"
$ws.Range("E17").Value = "  -3.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.772"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.103.76"
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.35"
$ws.Range("E22").Value = "  -2.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.64"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.04"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.92"
$ws.Range("E25").Value = "  -4.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("E27").Value = "  -4.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.78"
$ws.Range("E28").Value = "  -6.34%  "

$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.58"
$ws.Range("E30").Value = "  +4.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.65"
$ws.Range("E31").Value = "  -4.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.03"
$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  -4.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -3.59%  "

$ws.Range("E36").Value = "  -5.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("E37").Value = "  -5.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.14"
$ws.Range("E38").Value = "  -7.63%  "

$ws.Range("E39").Value = "  -5.52%  "

$ws.Range("E40").Value = "  -3.67%  "

$ws.Range("E41").Value = "  -3.61%  "

$ws.Range("E42").Value = "  -8.26%  "

$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.958.53"
$ws.Range("E44").Value = "  -2.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -2.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.30"
$ws.Range("E46").Value = "  -7.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.51"
$ws.Range("E47").Value = "  -6.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -5.50%  "

$ws.Range("E49").Value = "  -3.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.489.52"
$ws.Range("E50").Value = "  -2.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.55"
$ws.Range("E51").Value = "  -7.43%  "
